$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.703.46"
$ws.Range("D3").Value = "3.847.41"
$ws.Range("E3").Value = "  +2.94%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.96%  "
$ws.Range("D7").Value = "3.846.55"
$ws.Range("E7").Value = "  +2.95%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.84%  "
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").Value = "4.492.85"
$ws.Range("E15").Value = "  +2.99%  "
$ws.Range("D16").Value = "3.842.18"
$ws.Range("E16").Value = "  +3.47%  "
$ws.Range("D17").Value = "68.850.32"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.71%  "
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "486.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("E24").Value = "  +6.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("E26").Value = "  -1.64%  "
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.88%  "
$ws.Range("D32").Value = "3.997.75"
$ws.Range("E32").Value = "  +2.96%  "
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "32.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.80%  "
$ws.Range("D35").Value = "3.792.34"
$ws.Range("E35").Value = "  +3.32%  "
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("E37").Value = "  +1.62%  "
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("E39").Value = "  -0.63%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  -1.96%  "
$ws.Range("E42").Value = "  -1.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "433.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.41%  "
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.841.64"
$ws.Range("E48").Value = "  +2.26%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +12.88%  "
